{"js": "// Update the three-digit-by-one-digit multiplication problems with new\n// operands/results. Each original expression is unique in the document,\n// so a scoped find-and-replace (matchCase, exact text) for each pair is\n// sufficient and safe.\nconst replacements = [\n  [\"284\u00d79=\", \"327\u00d74=\"],\n  [\"771\u00d78=\", \"357\u00d74=\"],\n  [\"694\u00d73=\", \"988\u00d74=\"],\n  [\"339\u00d76=\", \"142\u00d72=\"],\n  [\"590\u00d79=\", \"578\u00d78=\"],\n  [\"700\u00d76=\", \"753\u00d75=\"],\n  [\"794\u00d74=\", \"131\u00d76=\"],\n  [\"495\u00d79=\", \"214\u00d75=\"],\n  [\"453\u00d76=\", \"249\u00d75=\"],\n  [\"710\u00d79=\", \"798\u00d74=\"],\n  [\"677\u00d74=\", \"470\u00d77=\"],\n  [\"219\u00d75=\", \"867\u00d76=\"],\n  [\"284\u00d72=\", \"954\u00d74=\"],\n  [\"846\u00d73=\", \"240\u00d77=\"],\n  [\"404\u00d72=\", \"389\u00d73=\"],\n  [\"638\u00d77=\", \"342\u00d77=\"],\n  [\"768\u00d78=\", \"452\u00d72=\"],\n  [\"583\u00d75=\", \"281\u00d74=\"],\n  [\"167\u00d79=\", \"467\u00d77=\"],\n  [\"153\u00d72=\", \"440\u00d73=\"],\n  [\"228\u00d75=\", \"668\u00d74=\"],\n  [\"583\u00d78=\", \"512\u00d76=\"],\n  [\"674\u00d73=\", \"677\u00d79=\"],\n  [\"758\u00d79=\", \"778\u00d78=\"],\n  [\"522\u00d73=\", \"806\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the three-digit-by-one-digit multiplication problems with new\n# operands/results. Each original expression is unique in the document,\n# so a scoped Find/Replace (exact match, match case) for each pair is\n# sufficient and safe.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"284\u00d79=\", \"327\u00d74=\"),\n    @(\"771\u00d78=\", \"357\u00d74=\"),\n    @(\"694\u00d73=\", \"988\u00d74=\"),\n    @(\"339\u00d76=\", \"142\u00d72=\"),\n    @(\"590\u00d79=\", \"578\u00d78=\"),\n    @(\"700\u00d76=\", \"753\u00d75=\"),\n    @(\"794\u00d74=\", \"131\u00d76=\"),\n    @(\"495\u00d79=\", \"214\u00d75=\"),\n    @(\"453\u00d76=\", \"249\u00d75=\"),\n    @(\"710\u00d79=\", \"798\u00d74=\"),\n    @(\"677\u00d74=\", \"470\u00d77=\"),\n    @(\"219\u00d75=\", \"867\u00d76=\"),\n    @(\"284\u00d72=\", \"954\u00d74=\"),\n    @(\"846\u00d73=\", \"240\u00d77=\"),\n    @(\"404\u00d72=\", \"389\u00d73=\"),\n    @(\"638\u00d77=\", \"342\u00d77=\"),\n    @(\"768\u00d78=\", \"452\u00d72=\"),\n    @(\"583\u00d75=\", \"281\u00d74=\"),\n    @(\"167\u00d79=\", \"467\u00d77=\"),\n    @(\"153\u00d72=\", \"440\u00d73=\"),\n    @(\"228\u00d75=\", \"668\u00d74=\"),\n    @(\"583\u00d78=\", \"512\u00d76=\"),\n    @(\"674\u00d73=\", \"677\u00d79=\"),\n    @(\"758\u00d79=\", \"778\u00d78=\"),\n    @(\"522\u00d73=\", \"806\u00d77=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
